$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.057237514097443
$ws.Range("D2").Value = 1.063801908004986
$ws.Range("E2").Value = 1.063362337245781
$ws.Range("F2").Value = 1.073856695906644
$ws.Range("I2").Value = 1.045622341484798
$ws.Range("J2").Value = 1.062234981983867
$ws.Range("K2").Value = 1.066519563723911
$ws.Range("L2").Value = 1.066081181700343
$ws.Range("M2").Value = 1.076547447634795
$ws.Range("N2").Value = 1.024472112371567
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.058389387893265
$ws.Range("D3").Value = 1.064831279183767
$ws.Range("E3").Value = 1.064378952264376
$ws.Range("F3").Value = 1.074940615255715
$ws.Range("I3").Value = 1.045867685861676
$ws.Range("J3").Value = 1.063038272645739
$ws.Range("K3").Value = 1.067363478996863
$ws.Range("L3").Value = 1.066912286670341
$ws.Range("M3").Value = 1.077447725449697
$ws.Range("N3").Value = 1.024741979345802
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.059134532509515
$ws.Range("D4").Value = 1.065497445981669
$ws.Range("E4").Value = 1.065036930321282
$ws.Range("F4").Value = 1.075642133877135
$ws.Range("I4").Value = 1.046024458939281
$ws.Range("J4").Value = 1.063557293716238
$ws.Range("K4").Value = 1.067909034190959
$ws.Range("L4").Value = 1.067449615381198
$ws.Range("M4").Value = 1.07802981516089
$ws.Range("N4").Value = 1.024916265018928
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.059447745538105
$ws.Range("D5").Value = 1.065777525366155
$ws.Range("E5").Value = 1.065313582761335
$ws.Range("F5").Value = 1.075937088474755
$ws.Range("I5").Value = 1.046089891943141
$ws.Range("J5").Value = 1.063775308047259
$ws.Range("K5").Value = 1.068138262631847
$ws.Range("L5").Value = 1.067675400381879
$ws.Range("M5").Value = 1.078274418055363
$ws.Range("N5").Value = 1.024989454215866
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.059500332724777
$ws.Range("D6").Value = 1.06582455327257
$ws.Range("E6").Value = 1.065360036182406
$ws.Range("F6").Value = 1.075986614775712
$ws.Range("I6").Value = 1.046100850610201
$ws.Range("J6").Value = 1.063811902939197
$ws.Range("K6").Value = 1.068176743890559
$ws.Range("L6").Value = 1.067713304349219
$ws.Range("M6").Value = 1.078315481651501
$ws.Range("N6").Value = 1.025001738283833
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.059138717856104
$ws.Range("D7").Value = 1.065501188326083
$ws.Range("E7").Value = 1.065040626813443
$ws.Range("F7").Value = 1.075646074932568
$ws.Range("I7").Value = 1.04602533512249
$ws.Range("J7").Value = 1.063560207549765
$ws.Range("K7").Value = 1.067912097635539
$ws.Range("L7").Value = 1.067452632755555
$ws.Range("M7").Value = 1.078033083979415
$ws.Range("N7").Value = 1.024917243292858
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.057626836081641
$ws.Range("D8").Value = 1.064149769283717
$ws.Range("E8").Value = 1.06370587396442
$ws.Range("F8").Value = 1.074222980511772
$ws.Range("I8").Value = 1.045705666910257
$ws.Range("J8").Value = 1.062506615677772
$ws.Range("K8").Value = 1.066804875471917
$ws.Range("L8").Value = 1.066362151162884
$ws.Range("M8").Value = 1.076851794064369
$ws.Range("N8").Value = 1.024563384821093
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.054961176629299
$ws.Range("D9").Value = 1.061769114183996
$ws.Range("E9").Value = 1.061355089205958
$ws.Range("F9").Value = 1.071716448565534
$ws.Range("I9").Value = 1.045127201720915
$ws.Range("J9").Value = 1.060644207526359
$ws.Range("K9").Value = 1.064849859400318
$ws.Range("L9").Value = 1.064437114892364
$ws.Range("M9").Value = 1.07476675369615
$ws.Range("N9").Value = 1.023937264458154
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.053182983392495
$ws.Range("D10").Value = 1.060182478382355
$ws.Range("E10").Value = 1.059788707427233
$ws.Range("F10").Value = 1.070046178876109
$ws.Range("I10").Value = 1.04473136322312
$ws.Range("J10").Value = 1.059398654043958
$ws.Range("K10").Value = 1.0635438439393
$ws.Range("L10").Value = 1.063151410720979
$ws.Range("M10").Value = 1.073374387683118
$ws.Range("N10").Value = 1.023518116259758
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.052412734387046
$ws.Range("D11").Value = 1.05949555400821
$ws.Range("E11").Value = 1.059110634633653
$ws.Range("F11").Value = 1.069323107136654
$ws.Range("I11").Value = 1.044557542794952
$ws.Range("J11").Value = 1.058858374129905
$ws.Range("K11").Value = 1.062977686600033
$ws.Range("L11").Value = 1.062594125777811
$ws.Range("M11").Value = 1.072770918613993
$ws.Range("N11").Value = 1.023336208158461
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.052126586286614
$ws.Range("D12").Value = 1.059240414138964
$ws.Range("E12").Value = 1.058858794928248
$ws.Range("F12").Value = 1.069054550508177
$ws.Range("I12").Value = 1.044492614528512
$ws.Range("J12").Value = 1.058657547279299
$ws.Range("K12").Value = 1.062767293287819
$ws.Range("L12").Value = 1.062387039741583
$ws.Range("M12").Value = 1.07254667788766
$ws.Range("N12").Value = 1.023268577048477
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.052187968007258
$ws.Range("D13").Value = 1.059295141862269
$ws.Range("E13").Value = 1.058912814198051
$ws.Range("F13").Value = 1.069112155720725
$ws.Range("I13").Value = 1.044506558317995
$ws.Range("J13").Value = 1.058700631792995
$ws.Range("K13").Value = 1.06281242777735
$ws.Range("L13").Value = 1.062431464280418
$ws.Range("M13").Value = 1.072594782159316
$ws.Range("N13").Value = 1.023283086974161
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.05238908218725
$ws.Range("D14").Value = 1.059474463775634
$ws.Range("E14").Value = 1.059089816942564
$ws.Range("F14").Value = 1.069300907667145
$ws.Range("I14").Value = 1.04455218322242
$ws.Range("J14").Value = 1.058841776636791
$ws.Range("K14").Value = 1.062960297410434
$ws.Range("L14").Value = 1.062577009730356
$ws.Range("M14").Value = 1.072752384551566
$ws.Range("N14").Value = 1.023330619022255
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.052512989490719
$ws.Range("D15").Value = 1.059584951821538
$ws.Range("E15").Value = 1.059198877702547
$ws.Range("F15").Value = 1.069417207189319
$ws.Range("I15").Value = 1.044580246046155
$ws.Range("J15").Value = 1.0589287216898
$ws.Range("K15").Value = 1.063051391867317
$ws.Range("L15").Value = 1.06266667372909
$ws.Range("M15").Value = 1.072849477249647
$ws.Range("N15").Value = 1.023359896821392
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.053234096279044
$ws.Range("D16").Value = 1.060228069342016
$ws.Range("E16").Value = 1.059833712657498
$ws.Range("F16").Value = 1.070094170218979
$ws.Range("I16").Value = 1.044742848124455
$ws.Range("J16").Value = 1.059434490642463
$ws.Range("K16").Value = 1.063581404324991
$ws.Range("L16").Value = 1.063188383910094
$ws.Range("M16").Value = 1.073414425982969
$ws.Range("N16").Value = 1.023530180165328
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.053686351731312
$ws.Range("D17").Value = 1.060631506229646
$ws.Range("E17").Value = 1.060231975957163
$ws.Range("F17").Value = 1.070518855323767
$ws.Range("I17").Value = 1.044844196222591
$ws.Range("J17").Value = 1.05975149220772
$ws.Range("K17").Value = 1.063913694354835
$ws.Range("L17").Value = 1.063515486906575
$ws.Range("M17").Value = 1.073768651674345
$ws.Range("N17").Value = 1.023636883441805
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.053950117796511
$ws.Range("D18").Value = 1.06086683376537
$ws.Range("E18").Value = 1.060464293682082
$ws.Range("F18").Value = 1.070766583019036
$ws.Range("I18").Value = 1.044903077388844
$ws.Range("J18").Value = 1.059936302519672
$ws.Range("K18").Value = 1.06410745140279
$ws.Range("L18").Value = 1.063706225918551
$ws.Range("M18").Value = 1.07397521100843
$ws.Range("N18").Value = 1.023699081704789
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.054040050687642
$ws.Range("D19").Value = 1.060947076021412
$ws.Range("E19").Value = 1.060543511001181
$ws.Range("F19").Value = 1.070851054501352
$ws.Range("I19").Value = 1.044923114747606
$ws.Range("J19").Value = 1.059999302584356
$ws.Range("K19").Value = 1.064173507030428
$ws.Range("L19").Value = 1.063771253740418
$ws.Range("M19").Value = 1.074045633171166
$ws.Range("N19").Value = 1.023720282931967
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.053637831791868
$ws.Range("D20").Value = 1.060588220277616
$ws.Range("E20").Value = 1.060189244269235
$ws.Range("F20").Value = 1.070473288973862
$ws.Range("I20").Value = 1.044833346679647
$ws.Range("J20").Value = 1.059717490391158
$ws.Range("K20").Value = 1.063878049206674
$ws.Range("L20").Value = 1.063480397504818
$ws.Range("M20").Value = 1.073730652252571
$ws.Range("N20").Value = 1.023625439323381
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.052329860280459
$ws.Range("D21").Value = 1.059421657547417
$ws.Range("E21").Value = 1.059037693311988
$ws.Range("F21").Value = 1.069245324253953
$ws.Range("I21").Value = 1.044538757863605
$ws.Range("J21").Value = 1.058800216945642
$ws.Range("K21").Value = 1.062916756178256
$ws.Range("L21").Value = 1.062534152589683
$ws.Range("M21").Value = 1.072705976937923
$ws.Range("N21").Value = 1.023316623740691
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.051507235359358
$ws.Range("D22").Value = 1.058688276482649
$ws.Range("E22").Value = 1.058313821893855
$ws.Range("F22").Value = 1.068473395055661
$ws.Range("I22").Value = 1.044351434345269
$ws.Range("J22").Value = 1.058222664048588
$ws.Range("K22").Value = 1.062311790076006
$ws.Range("L22").Value = 1.06193871538683
$ws.Range("M22").Value = 1.072061228289563
$ws.Range("N22").Value = 1.023122098377084
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.051943348658461
$ws.Range("D23").Value = 1.059077047942968
$ws.Range("E23").Value = 1.058697545348018
$ws.Range("F23").Value = 1.068882596149521
$ws.Range("I23").Value = 1.044450937534534
$ws.Range("J23").Value = 1.058528914312804
$ws.Range("K23").Value = 1.062632547593445
$ws.Range("L23").Value = 1.062254415039878
$ws.Range("M23").Value = 1.072403068799269
$ws.Range("N23").Value = 1.023225254187833
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.053659755946644
$ws.Range("D24").Value = 1.060607779306123
$ws.Range("E24").Value = 1.060208552825394
$ws.Range("F24").Value = 1.070493878396117
$ws.Range("I24").Value = 1.044838249842348
$ws.Range("J24").Value = 1.059732854630436
$ws.Range("K24").Value = 1.063894155906448
$ws.Range("L24").Value = 1.063496253063634
$ws.Range("M24").Value = 1.073747822723807
$ws.Range("N24").Value = 1.023630610551264
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.055650500810026
$ws.Range("D25").Value = 1.062384487038294
$ws.Range("E25").Value = 1.06196268005027
$ws.Range("F25").Value = 1.072364313448021
$ws.Range("I25").Value = 1.045278546030333
$ws.Range("J25").Value = 1.061126379166458
$ws.Range("K25").Value = 1.06535574755417
$ws.Range("L25").Value = 1.064935195420234
$ws.Range("M25").Value = 1.075306197757473
$ws.Range("N25").Value = 1.024099437140107

Write-Host "Updated vm_pu values for 380 kV case"
